$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column A ("Employee_ID") in front of the existing data,
# shifting Name/Email/Department/Delegate Comments/Predicted Category
# one column to the right.
$ws.Columns("A").Insert()

# Drop the Alice Brown and Charlie Wilson rows (original rows 3 and 4,
# i.e. rows 3:4 after the column insert above - row numbers are unaffected
# by a column insert).
$ws.Rows("3:4").Delete()

# Header for the new column.
$ws.Range("A1").Value = "Employee_ID"
$ws.Range("A1").Font.Bold = $true

# Employee id value for the remaining (Jane Smith) row - must stay text
# ("002"), not be coerced into the number 2. Briefly mark the cell as
# Text, assign the literal value, then restore the default number
# format so the cell keeps the workbook's ordinary (unstyled) look.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "002"
$ws.Range("A2").NumberFormat = ""

Write-Output "done"
